$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D='27.163.66'; E='  -1.96%  '},
    @{Row=3; D='1.818.18'; E='  -1.84%  '},
    @{Row=4; D='1.007'; E='  -0.93%  '},
    @{Row=5; D='314.01'; E='  -1.77%  '},
    @{Row=6; E='  -1.02%  '},
    @{Row=7; D='0.4267'; E='  -2.39%  '},
    @{Row=8; D='0.3667'; E='  -2.79%  '},
    @{Row=9; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='45.67'; E='  -2.04%  '},
    @{Row=10; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.07221'; E='  -2.54%  '},
    @{Row=11; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.8599'; E='  -2.43%  '},
    @{Row=12; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='20.95'; E='  -2.61%  '},
    @{Row=13; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.829.23'; E='  -1.50%  '},
    @{Row=14; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='6.660'; E='  -1.54%  '},
    @{Row=15; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07094'; E='  +0.01%  '},
    @{Row=16; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.301'; E='  -3.02%  '},
    @{Row=17; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='88.09'; E='  -0.13%  '},
    @{Row=18; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.008'; E='  -1.12%  '},
    @{Row=19; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000008847'; E='  -1.62%  '},
    @{Row=20; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.005'; E='  -1.02%  '},
    @{Row=21; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='15.04'; E='  -2.24%  '},
    @{Row=22; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='27.222.69'; E='  -1.79%  '},
    @{Row=23; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='5.138'; E='  -2.24%  '},
    @{Row=24; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='10.85'; E='  -2.41%  '},
    @{Row=25; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.050.86'; E='  -1.68%  '},
    @{Row=26; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='2.003'; E='  -1.07%  '},
    @{Row=27; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='153.54'; E='  -1.92%  '},
    @{Row=28; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='18.34'; E='  -1.43%  '},
    @{Row=29; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='2.109'; E='  +6.37%  '},
    @{Row=30; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='5.219'; E='  -3.92%  '},
    @{Row=31; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='116.06'; E='  -3.66%  '},
    @{Row=32; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.08895'; E='  -1.56%  '},
    @{Row=33; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.193'; E='  -2.72%  '},
    @{Row=34; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.7562'; E='  -1.54%  '},
    @{Row=35; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='4.453'; E='  -2.09%  '},
    @{Row=36; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.829'; E='  -5.56%  '},
    @{Row=37; B='Frax'; C='https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D='1.005'; E='  -1.15%  '},
    @{Row=38; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.113'; E='  -1.99%  '},
    @{Row=39; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.01968'; E='  -0.32%  '},
    @{Row=40; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.05255'; E='  -0.80%  '},
    @{Row=41; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.900'; E='  +1.17%  '},
    @{Row=42; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='7.142'; E='  +2.83%  '},
    @{Row=43; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1684'; E='  +0.56%  '},
    @{Row=44; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.5042'; E='  -2.62%  '},
    @{Row=45; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='8.584'; E='  -1.09%  '},
    @{Row=46; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='10.61'; E='  -0.81%  '},
    @{Row=47; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='106.49'; E='  -3.04%  '},
    @{Row=48; B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.4739'; E='  +0.30%  '},
    @{Row=49; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='1.005'; E='  -1.11%  '},
    @{Row=50; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.06386'; E='  -1.20%  '},
    @{Row=51; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='1.656'; E='  -2.97%  '}
)

foreach ($item in $data) {
    if ($item.ContainsKey("B")) { $ws.Cells.Item($item.Row, 2).Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Cells.Item($item.Row, 3).Value = $item.C }
    if ($item.ContainsKey("D")) {
        $cell = $ws.Cells.Item($item.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.Style = "Normal"
    }
    if ($item.ContainsKey("E")) { $ws.Cells.Item($item.Row, 5).Value = $item.E }
}
